$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B24").Value = "<50 Production`n<50 Services`n<25 Wholesale`n<15 Retail"
$ws.Range("C24").Value = "<50Millionlion bath Production, `n<50Millionlion bath Services, `n<50Millionlion bath Wholesale, `n<30Millionlion bath Retail"

$ws.Range("B25").Value = "≤51 Production<200, `n≤51 Services<200, `n≤26 Wholesale<200, `n≤16 Retail <30"
$ws.Range("C25").Value = "≤50 Production<200Millionlion bath, `n≤50 Services<200Millionlion bath, `n≤50 Wholesale<100Millionlion bath, `n≤30 Retail <60Millionlion bath"

$ws.Range("B26").Value = ">=200 Production, `n>=200 Services, `n>=200 Wholesale, `n>=30 Retail"
$ws.Range("C26").Value = ">=200Millionlion bath Production, `n>=200Millionlion bath Services, `n>=100Millionlion bath Wholesale, `n>=60Millionlion bath Retail"
